$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.307.75'
$ws.Range("E2").Value = '  +1.29%  '

$ws.Range("D3").Value = '1.833.49'
$ws.Range("E3").Value = '  +0.83%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.011'
$ws.Range("E4").Value = '  +0.86%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.74'
$ws.Range("E5").Value = '  +1.57%  '

$ws.Range("E6").Value = '  +0.74%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4735'
$ws.Range("E7").Value = '  +1.74%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3687'
$ws.Range("E8").Value = '  +0.70%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07448'
$ws.Range("E9").Value = '  +1.28%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8864'
$ws.Range("E10").Value = '  +1.69%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.50'
$ws.Range("E11").Value = '  +1.16%  '

$ws.Range("D12").Value = '1.883.48'
$ws.Range("E12").Value = '  +3.44%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07331'
$ws.Range("E13").Value = '  +3.10%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.430'
$ws.Range("E14").Value = '  +0.61%  '

$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '94.13'
$ws.Range("E15").Value = '  +2.94%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.567'
$ws.Range("E16").Value = '  +0.95%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.010'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008798'
$ws.Range("E18").Value = '  +1.16%  '

$ws.Range("E19").Value = '  +0.81%  '

$ws.Range("D20").Value = '27.585.43'
$ws.Range("E20").Value = '  +2.21%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.79'
$ws.Range("E21").Value = '  +0.90%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.288'
$ws.Range("E22").Value = '  -0.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.68'
$ws.Range("E23").Value = '  +0.98%  '

$ws.Range("D24").Value = '2.099.79'
$ws.Range("E24").Value = '  +2.54%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.896'
$ws.Range("E25").Value = '  +0.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.87'
$ws.Range("E26").Value = '  +0.76%  '

$ws.Range("E27").Value = '  +1.17%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.143'
$ws.Range("E28").Value = '  +0.21%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.233'
$ws.Range("E29").Value = '  -0.14%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.33'
$ws.Range("E30").Value = '  +0.83%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08994'
$ws.Range("E31").Value = '  +1.08%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7503'
$ws.Range("E32").Value = '  -1.11%  '

$ws.Range("E33").Value = '  +0.93%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.548'
$ws.Range("E34").Value = '  +0.96%  '

$ws.Range("E35").Value = '  +1.51%  '

$ws.Range("E36").Value = '  +0.89%  '

$ws.Range("E37").Value = '  +0.34%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05340'
$ws.Range("E38").Value = '  +0.84%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01957'
$ws.Range("E39").Value = '  +0.45%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.434'
$ws.Range("E40").Value = '  +3.45%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.964'
$ws.Range("E41").Value = '  -0.16%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.249'
$ws.Range("E42").Value = '  +1.06%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5293'
$ws.Range("E43").Value = '  +0.25%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1660'
$ws.Range("E44").Value = '  +0.12%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.495'
$ws.Range("E45").Value = '  +0.58%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4926'
$ws.Range("E46").Value = '  +1.31%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.54'
$ws.Range("E47").Value = '  +0.83%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '105.26'
$ws.Range("E48").Value = '  +1.77%  '

$ws.Range("E49").Value = '  +0.83%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.667'
$ws.Range("E50").Value = '  +0.06%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06299'
$ws.Range("E51").Value = '  +0.14%  '
